# Re-applies the refreshed "cryptos" price/volume snapshot.
# Values are prefixed with a leading apostrophe (exactly as typing into
# the Excel UI would) so numeric-looking strings such as "19.35",
# "1.00" or "0.0850" are stored as literal text instead of being
# auto-converted to numbers (which would drop formatting/trailing zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.281.14"
$ws.Range("E2").Value = "'  +0.27%  "
$ws.Range("D3").Value = "'1.590.05"
$ws.Range("E3").Value = "'  +0.45%  "
$ws.Range("E4").Value = "'  -0.21%  "
$ws.Range("D5").Value = "'213.13"
$ws.Range("E5").Value = "'  +1.61%  "
$ws.Range("E6").Value = "'  +0.59%  "
$ws.Range("E7").Value = "'  -0.20%  "
$ws.Range("E8").Value = "'  -0.05%  "
$ws.Range("E9").Value = "'  -0.27%  "
$ws.Range("D10").Value = "'19.35"
$ws.Range("E10").Value = "'  -0.91%  "
$ws.Range("D11").Value = "'0.0850"
$ws.Range("E11").Value = "'  +0.45%  "
$ws.Range("D12").Value = "'1.812.80"
$ws.Range("E12").Value = "'  +0.41%  "
$ws.Range("D13").Value = "'1.595.49"
$ws.Range("E13").Value = "'  +0.03%  "
$ws.Range("E14").Value = "'  -0.36%  "
$ws.Range("E15").Value = "'  +1.10%  "
$ws.Range("D16").Value = "'64.46"
$ws.Range("E16").Value = "'  -0.10%  "
$ws.Range("D17").Value = "'26.282.32"
$ws.Range("E17").Value = "'  +0.26%  "
$ws.Range("D18").Value = "'0.0₃0727"
$ws.Range("E18").Value = "'  -0.99%  "
$ws.Range("D19").Value = "'7.48"
$ws.Range("E19").Value = "'  +2.54%  "
$ws.Range("D20").Value = "'213.78"
$ws.Range("E20").Value = "'  +3.14%  "
$ws.Range("E21").Value = "'  -0.15%  "
$ws.Range("E22").Value = "'  +0.76%  "
$ws.Range("D23").Value = "'2.16"
$ws.Range("E23").Value = "'  -2.13%  "
$ws.Range("D24").Value = "'8.97"
$ws.Range("E24").Value = "'  +0.73%  "
$ws.Range("D25").Value = "'145.09"
$ws.Range("E25").Value = "'  +0.42%  "
$ws.Range("D27").Value = "'7.05"
$ws.Range("E27").Value = "'  +0.58%  "
$ws.Range("E28").Value = "'  -0.55%  "
$ws.Range("E29").Value = "'  -0.27%  "
$ws.Range("E30").Value = "'  -0.79%  "
$ws.Range("E31").Value = "'  +1.20%  "
$ws.Range("E32").Value = "'  -0.21%  "
$ws.Range("D33").Value = "'2.96"
$ws.Range("E33").Value = "'  +0.24%  "
$ws.Range("D34").Value = "'1.341.18"
$ws.Range("E34").Value = "'  +4.99%  "
$ws.Range("E35").Value = "'  -0.91%  "
$ws.Range("D37").Value = "'0.596"
$ws.Range("E37").Value = "'  -2.76%  "
$ws.Range("E38").Value = "'  -0.35%  "
$ws.Range("D39").Value = "'0.817"
$ws.Range("E39").Value = "'  -0.12%  "
$ws.Range("D40").Value = "'5.79"
$ws.Range("E40").Value = "'  +4.29%  "
$ws.Range("E41").Value = "'  -0.18%  "
$ws.Range("E42").Value = "'  -3.39%  "
$ws.Range("E43").Value = "'  +0.31%  "
$ws.Range("E44").Value = "'  -0.52%  "
$ws.Range("B45").Value = "'RocketPoolETH"
$ws.Range("C45").Value = "'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "'1.724.39"
$ws.Range("E45").Value = "'  +0.24%  "
$ws.Range("B46").Value = "'Aave"
$ws.Range("C46").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'61.79"
$ws.Range("E46").Value = "'  -0.89%  "
$ws.Range("D47").Value = "'87.32"
$ws.Range("E47").Value = "'  -2.13%  "
$ws.Range("D48").Value = "'1.50"
$ws.Range("E48").Value = "'  -3.96%  "
$ws.Range("E49").Value = "'  -0.53%  "
$ws.Range("D50").Value = "'0.0980"
$ws.Range("E50").Value = "'  -2.44%  "
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "'  -0.32%  "
